# Generate Report for Handoff
# - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview sheet (zh-cn + de-de columns) and on each language sheet's
#   Status column.
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   that were tied to that status also advance.
# - Because the new status text is shorter, the Status-holding columns that
#   were auto-fit to the old text shrink to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus          # zh-cn status
$overview.Range("F2").Value = $newStatus          # de-de status
$overview.Range("G2").Value = "2016-08-28 18:57:29"  # Latest HO Xliff Generate Date

# Columns E and F were sized to fit the old, longer status text; refit them
# to (approximately) the width that the shorter text now needs.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus              # Status
$zhcn.Range("H2").Value = "2016-08-28 18:57:25"   # Latest Handoff Datetime
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus              # Status
$dede.Range("H2").Value = "2016-08-28 18:57:29"   # Latest Handoff Datetime
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
